# C5-PowerPoint.pptx edit:
#  1. Slide 6's table switches to a different built-in table style.
#  2. The presentation's design colours switch from the "Integral" theme
#     palette over to the standard "Office Theme" palette (the slide
#     master / design's theme colours, i.e. ppt/theme/theme2.xml).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{11E4FAC3-CF12-439F-BC7F-D50654920598}")
    }
}

# --- 2. Swap the design's theme colours to the Office Theme palette ------
$scheme = $p.SlideMaster.ColorScheme

$scheme.Colors(1).RGB  = 0          # dk1      000000
$scheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2      44546A
$scheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$scheme.Colors(10).RGB = 4697456    # accent6  70AD47
$scheme.Colors(11).RGB = 12673797   # hlink    0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
